$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 336.66666
$ws.Range("I12").Value = 254
$ws.Range("J12").Value = 750
$ws.Range("K12").Value = 254
$ws.Range("L12").Value = 750
$ws.Range("M12").Value = -84
$ws.Range("N12").Value = -1090

$ws.Range("H32").Value = 3749.6667
$ws.Range("J32").Value = 3749.6667
$ws.Range("L32").Value = 3749.6667
$ws.Range("N32").Value = -4401.6667

$ws.Range("H33").Value = 318
$ws.Range("I33").Value = 266
$ws.Range("J33").Value = 508.66666
$ws.Range("K33").Value = 266
$ws.Range("L33").Value = 508.66666
$ws.Range("M33").Value = -37
$ws.Range("N33").Value = -966.66666

$ws.Range("H34").Value = 2795.111
$ws.Range("I34").Value = 2795.111
$ws.Range("K34").Value = 2795.111
$ws.Range("M34").Value = -2592.111

$ws.Range("H36").Value = 2795.111
$ws.Range("I36").Value = 2795.111
$ws.Range("K36").Value = 2795.111
$ws.Range("M36").Value = -2080.111

$ws.Range("H69").Value = 15110.889

$ws.Range("H72").Value = 15110.889

$ws.Range("H86").Value = 7937.25
$ws.Range("J86").Value = 8642.571
$ws.Range("L86").Value = 8642.571
$ws.Range("N86").Value = -10888.571

$ws.Range("H87").Value = 59987.5
$ws.Range("J87").Value = 114975
$ws.Range("L87").Value = 114975
$ws.Range("N87").Value = -117471

$ws.Range("H89").Value = 7937.25
$ws.Range("J89").Value = 8642.571
$ws.Range("L89").Value = 43212.855
$ws.Range("N89").Value = -54444.855

$ws.Range("H90").Value = 59987.5
$ws.Range("J90").Value = 114975
$ws.Range("L90").Value = 344925
$ws.Range("N90").Value = -357405

$ws.Range("H132").Value = 5653.143
$ws.Range("I132").Value = 2363.3076
$ws.Range("J132").Value = 10999.125
$ws.Range("K132").Value = 7089.9228
$ws.Range("L132").Value = 32997.375
$ws.Range("M132").Value = -4559.9228
$ws.Range("N132").Value = -38057.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2402.8125
$ws.Range("I2").Value = 2370.5
$ws.Range("K2").Value = 2370.5
$ws.Range("M2").Value = -2257.5

$ws.Range("H61").Value = 35002816
$ws.Range("I61").Value = 40003380
$ws.Range("K61").Value = 40003380
$ws.Range("M61").Value = -40003168

$ws.Range("H116").Value = 2402.8125
$ws.Range("I116").Value = 2370.5
$ws.Range("K116").Value = 2370.5
$ws.Range("M116").Value = -76.5

$ws.Range("I132").Value = 5747.65
$ws.Range("K132").Value = 17242.95
$ws.Range("M132").Value = -14712.95

$ws.Range("H136").Value = 35002816
$ws.Range("I136").Value = 40003380
$ws.Range("K136").Value = 120010140
$ws.Range("M136").Value = -120007590

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2402.8125
$ws.Range("I3").Value = 2370.5
$ws.Range("K3").Value = 2370.5
$ws.Range("M3").Value = -2256.5

$ws.Range("H134").Value = 5002729.5
$ws.Range("I134").Value = 1899.25
$ws.Range("K134").Value = 5697.75
$ws.Range("M134").Value = -3162.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43925604
$ws.Range("I31").Value = 52635590
$ws.Range("K31").Value = 52635590
$ws.Range("M31").Value = -52635295

$ws.Range("H34").Value = 43925604
$ws.Range("I34").Value = 52635590
$ws.Range("K34").Value = 52635590
$ws.Range("M34").Value = -52635388

$ws.Range("H94").Value = 1048.2778
$ws.Range("I94").Value = 1137.875
$ws.Range("J94").Value = 976.6
$ws.Range("K94").Value = 1137.875
$ws.Range("L94").Value = 976.6
$ws.Range("M94").Value = -686.875
$ws.Range("N94").Value = -1878.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 13657.391
$ws.Range("I134").Value = 1198.88
$ws.Range("K134").Value = 3596.64
$ws.Range("M134").Value = 1473.36

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3299.3333
$ws.Range("J80").Value = 3452.6667
$ws.Range("L80").Value = 3452.6667
$ws.Range("N80").Value = -5448.6667

$ws.Range("H83").Value = 3299.3333
$ws.Range("J83").Value = 3452.6667
$ws.Range("L83").Value = 17263.3335
$ws.Range("N83").Value = -27247.3335

$ws.Range("H132").Value = 7639438.5
$ws.Range("I132").Value = 3251.8823
$ws.Range("K132").Value = 9755.6469
$ws.Range("M132").Value = -7225.6469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4293.5557
$ws.Range("I22").Value = 3829.8
$ws.Range("J22").Value = 4873.25
$ws.Range("K22").Value = 3829.8
$ws.Range("L22").Value = 4873.25
$ws.Range("M22").Value = -3534.8
$ws.Range("N22").Value = -5463.25

$ws.Range("H27").Value = 4293.5557
$ws.Range("I27").Value = 3829.8
$ws.Range("J27").Value = 4873.25
$ws.Range("K27").Value = 3829.8
$ws.Range("L27").Value = 4873.25
$ws.Range("M27").Value = -3722.8
$ws.Range("N27").Value = -5087.25

$ws.Range("H40").Value = 3038.182
$ws.Range("I40").Value = 2692.05
$ws.Range("K40").Value = 2692.05
$ws.Range("M40").Value = -2556.05

$ws.Range("H82").Value = 7254.222
$ws.Range("I82").Value = 5125
$ws.Range("J82").Value = 8957.6
$ws.Range("K82").Value = 5125
$ws.Range("L82").Value = 8957.6
$ws.Range("M82").Value = -4764
$ws.Range("N82").Value = -9679.6

$ws.Range("H85").Value = 7254.222
$ws.Range("I85").Value = 5125
$ws.Range("J85").Value = 8957.6
$ws.Range("K85").Value = 5125
$ws.Range("L85").Value = 8957.6
$ws.Range("M85").Value = -3877
$ws.Range("N85").Value = -11453.6

$ws.Range("H93").Value = 6953194.5
$ws.Range("I93").Value = 3666.6667
$ws.Range("K93").Value = 3666.6667
$ws.Range("M93").Value = -2418.6667

$ws.Range("H100").Value = 19253938
$ws.Range("I100").Value = 3719.2856
$ws.Range("K100").Value = 3719.2856
$ws.Range("M100").Value = -3178.2856

$ws.Range("H122").Value = 3431.617
$ws.Range("I122").Value = 3286.3777
$ws.Range("K122").Value = 9859.133099999999
$ws.Range("M122").Value = -7409.133099999999

$ws.Range("H132").Value = 4537.75
$ws.Range("I132").Value = 2636.5715
$ws.Range("J132").Value = 7199.4
$ws.Range("K132").Value = 7909.7145
$ws.Range("L132").Value = 21598.2
$ws.Range("M132").Value = -5379.7145
$ws.Range("N132").Value = -26658.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1578.6471
$ws.Range("I81").Value = 1427.3125
$ws.Range("K81").Value = 2854.625
$ws.Range("M81").Value = -1793.625

$ws.Range("H84").Value = 1578.6471
$ws.Range("I84").Value = 1427.3125
$ws.Range("K84").Value = 14273.125
$ws.Range("M84").Value = -8969.125

$ws.Range("H126").Value = 11373.7
$ws.Range("I126").Value = 11551.333
$ws.Range("K126").Value = 34653.999
$ws.Range("M126").Value = -32183.999

$ws.Range("H132").Value = 771243.1
$ws.Range("I132").Value = 1616.1
$ws.Range("K132").Value = 4848.299999999999
$ws.Range("M132").Value = -2318.299999999999

$ws.Range("H136").Value = 401702.72
$ws.Range("I136").Value = 1720.1305
$ws.Range("K136").Value = 5160.3915
$ws.Range("M136").Value = -2610.3915
